$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Frequency -> hz / 50 ---
$ws.Range("A3").Copy()
$ws.Range("H3").PasteSpecial(-4104)
$ws.Range("H3").Value = "hz"

$ws.Range("A3").Copy()
$ws.Range("I3").PasteSpecial(-4104)
$ws.Range("I3").Value = 50

# --- Row 4: Temperature -> C / 25 ---
$ws.Range("A4").Copy()
$ws.Range("H4").PasteSpecial(-4104)
$ws.Range("H4").Value = "C"

$ws.Range("A4").Copy()
$ws.Range("I4").PasteSpecial(-4104)
$ws.Range("I4").Value = 25

# --- Row 5: Sound -> dB / 100 ---
$ws.Range("A5").Copy()
$ws.Range("H5").PasteSpecial(-4104)
$ws.Range("H5").Value = "dB"

$ws.Range("A5").Copy()
$ws.Range("I5").PasteSpecial(-4104)
$ws.Range("I5").Value = 100

# --- Update selection to match the author's final cursor position ---
$ws.Range("H11").Select()
